$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 2721
$ws.Range("F4").Value = 321
$ws.Range("F5").Value = 220
$ws.Range("F6").Value = 534
$ws.Range("F7").Value = 1312
$ws.Range("F8").Value = 622
$ws.Range("F9").Value = 335
$ws.Range("F11").Value = 148
$ws.Range("F12").Value = 399
$ws.Range("F13").Value = 6048
$ws.Range("F14").Value = 107
$ws.Range("F16").Value = 1874
$ws.Range("F17").Value = 4514
$ws.Range("F18").Value = 455
$ws.Range("F21").Value = 5263
$ws.Range("F22").Value = 6844
$ws.Range("F24").Value = 1075
$ws.Range("F25").Value = 733
$ws.Range("F26").Value = 3919
$ws.Range("F27").Value = 533
$ws.Range("F28").Value = 71
$ws.Range("F29").Value = 213
$ws.Range("F30").Value = 141
$ws.Range("F31").Value = 1031
$ws.Range("F32").Value = 1473
$ws.Range("F33").Value = 531
$ws.Range("F34").Value = 648
$ws.Range("F35").Value = 1656
$ws.Range("F36").Value = 226
$ws.Range("F37").Value = 1833
$ws.Range("F38").Value = 1205
$ws.Range("F39").Value = 1337
$ws.Range("F40").Value = 667
$ws.Range("F41").Value = 274
$ws.Range("F42").Value = 269
$ws.Range("F43").Value = 3593
$ws.Range("F44").Value = 150
$ws.Range("F45").Value = 329
$ws.Range("F46").Value = 434
$ws.Range("F48").Value = 79
$ws.Range("F49").Value = 3927

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 1253
$ws.Range("F5").Value = 43
$ws.Range("F7").Value = 5
$ws.Range("F9").Value = 20

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 4264

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 4264
$ws.Range("F3").Value = 2721
$ws.Range("F5").Value = 321
$ws.Range("F6").Value = 1253
$ws.Range("F8").Value = 220
$ws.Range("F9").Value = 534
$ws.Range("F11").Value = 1312
$ws.Range("F12").Value = 20
$ws.Range("F13").Value = 622
$ws.Range("F14").Value = 335
$ws.Range("F16").Value = 148
$ws.Range("F17").Value = 399
$ws.Range("F18").Value = 107
$ws.Range("F20").Value = 1874
$ws.Range("F21").Value = 4515
$ws.Range("F22").Value = 5263
$ws.Range("F23").Value = 5263
$ws.Range("F25").Value = 1075
$ws.Range("F26").Value = 733
$ws.Range("F27").Value = 3920
$ws.Range("F28").Value = 533
$ws.Range("F30").Value = 141
$ws.Range("F31").Value = 1031
$ws.Range("F32").Value = 1473
$ws.Range("F33").Value = 531
$ws.Range("F34").Value = 648
$ws.Range("F35").Value = 1656
$ws.Range("F36").Value = 1833
$ws.Range("F39").Value = 667
$ws.Range("F41").Value = 275
$ws.Range("F43").Value = 3593
$ws.Range("F45").Value = 150
$ws.Range("F46").Value = 329
$ws.Range("F47").Value = 434
$ws.Range("F48").Value = 79
$ws.Range("F50").Value = 3927

